# Commiting Customer + Account files
#
# Appends 41 new customer rows (rows 81-121) to Sheet0 of the
# "Individual Customers" workbook, extending the used range from
# A1:C80 to A1:C121. Columns:
#   A = Company        B = Customer_ID        C = PD (department code)
#
# All three columns in this sheet are authored as shared-string TEXT,
# even though the values look numeric (e.g. "17704760"). A plain
# `Range.Value = "17704760"` assignment would be auto-coerced to a
# number by Excel's normal type inference, and forcing text via
# NumberFormat="@" (or a leading apostrophe) stamps the cell with a new
# style index, neither of which matches the source diff. Instead, each
# value is staged as a `="text"` formula in a scratch area, then copied
# and pasted as values (PasteSpecial xlPasteValues = -4163) into the
# destination cell. A formula's string result pastes as genuine shared
# -string text without touching cell formatting/styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: (row, A-Company, B-Customer_ID, C-PD)
$rows = @(
  @(81,  "118500", "17704760", "6004"),
  @(82,  "118518", "17704761", "6020"),
  @(83,  "118498", "17704762", "1001"),
  @(84,  "118452", "17704763", "1001"),
  @(85,  "118518", "17704764", "6012"),
  @(86,  "118500", "17704767", "6004"),
  @(87,  "118518", "17704768", "6020"),
  @(88,  "118498", "17704769", "1001"),
  @(89,  "118500", "17704772", "6004"),
  @(90,  "118518", "17704773", "6020"),
  @(91,  "118498", "17704774", "1001"),
  @(92,  "118500", "17704775", "6004"),
  @(93,  "118518", "17704776", "6020"),
  @(94,  "118498", "17704777", "1001"),
  @(95,  "118518", "17704779", "6020"),
  @(96,  "118518", "17704781", "6020"),
  @(97,  "118518", "17704783", "6020"),
  @(98,  "118518", "17704784", "6020"),
  @(99,  "118518", "17704785", "6020"),
  @(100, "118518", "17704786", "6020"),
  @(101, "118518", "17704787", "6020"),
  @(102, "118518", "17704789", "6020"),
  @(103, "118500", "17704791", "6004"),
  @(104, "118498", "17704793", "1001"),
  @(105, "118452", "17704794", "1001"),
  @(106, "118518", "17704795", "6012"),
  @(107, "118500", "17704800", "6004"),
  @(108, "118518", "17704801", "6020"),
  @(109, "118498", "17704802", "1001"),
  @(110, "118452", "17704803", "1001"),
  @(111, "118518", "17704804", "6012"),
  @(112, "118500", "17704808", "6004"),
  @(113, "118518", "17704809", "6020"),
  @(114, "118498", "17704810", "1001"),
  @(115, "118452", "17704811", "1001"),
  @(116, "118518", "17704812", "6012"),
  @(117, "118448", "17704813", "1047"),
  @(118, "118518", "17704815", "1035"),
  @(119, "118452", "17704816", "1150"),
  @(120, "118448", "17704817", "1068"),
  @(121, "118448", "17704818", "1005")
)

$n = $rows.Length

# Scratch columns Z, AA, AB (26, 27, 28) - far away from the data table,
# one staging row per new data row. Each holds a `="value"` formula
# whose cached/evaluated result is text.
for ($i = 0; $i -lt $n; $i++) {
  $ws.Cells.Item($i + 1, 26).Formula = '="' + $rows[$i][1] + '"'
  $ws.Cells.Item($i + 1, 27).Formula = '="' + $rows[$i][2] + '"'
  $ws.Cells.Item($i + 1, 28).Formula = '="' + $rows[$i][3] + '"'
}

$firstRow = $rows[0][0]
$lastRow = $rows[$n - 1][0]

# Column A (Company)
$ws.Range($ws.Cells.Item(1, 26), $ws.Cells.Item($n, 26)).Copy()
$ws.Range($ws.Cells.Item($firstRow, 1), $ws.Cells.Item($lastRow, 1)).PasteSpecial(-4163)

# Column B (Customer_ID)
$ws.Range($ws.Cells.Item(1, 27), $ws.Cells.Item($n, 27)).Copy()
$ws.Range($ws.Cells.Item($firstRow, 2), $ws.Cells.Item($lastRow, 2)).PasteSpecial(-4163)

# Column C (PD)
$ws.Range($ws.Cells.Item(1, 28), $ws.Cells.Item($n, 28)).Copy()
$ws.Range($ws.Cells.Item($firstRow, 3), $ws.Cells.Item($lastRow, 3)).PasteSpecial(-4163)

# Clean up the scratch area so it doesn't linger in the saved workbook.
$ws.Range($ws.Cells.Item(1, 26), $ws.Cells.Item($n, 28)).Value = ""
